$d = $word.ActiveDocument

# --- Change 1: paragraph "Bei Mehrfachverwendung von Kanten..." ---
# Merge the three runs into a single run and apply strikethrough
# (to both the run and the paragraph mark formatting).
$target = "Bei Mehrfachverwendung von Kanten: Kreuzungsvermeidung ( und Sperre nach oben) verbessern!"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq $target) {
        $full = $p.Range
        $body = $d.Range($full.Start, $full.End - 1)
        $body.Delete()
        $ins = $d.Range($full.Start, $full.Start)
        $ins.InsertAfter($target)

        $p.Range.Font.StrikeThrough = $true
        $p.Range.Font.Size = 10
        $p.Range.LanguageID = "en-GB"
        break
    }
}

# --- Change 2: remove the paragraph "Nach gefundnen Positionen nacher Platzierung nachbessern?" ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Nach gefundnen Positionen nacher Platzierung nachbessern?") {
        $p.Range.Delete()
        break
    }
}
